$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Location")

# Rename the existing "Home.AddTank" key to "Home.AddTank.Add"
$ws.Range("A2").Value = "Home.AddTank.Add"

# Add three new rows describing the "Reset Data" modal, following the
# same alternating-style pattern as the existing table rows (odd rows
# use the style of row 19, even rows use the style of row 20).
$ws.Range("A19:G19").Copy()
$ws.Range("A21:G21").PasteSpecial(-4122)
$ws.Rows(21).RowHeight = 20.25

$ws.Range("A20:G20").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)
$ws.Rows(22).RowHeight = 20.25

$ws.Range("A19:G19").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)
$ws.Rows(23).RowHeight = 20.25

$ws.Range("A21").Value = "Home.AddTank.ResetData"
$ws.Range("B21").Value = "XPath"
$ws.Range("C21").Value = "//div[@id='ResetData']//a[text()='Reset Data']"

$ws.Range("A22").Value = "Home.AddTank.ResetModal"
$ws.Range("B22").Value = "XPath"
$ws.Range("C22").Value = "//div[@id='ResetData']"

$ws.Range("A23").Value = "Home.AddTank.CloseModal"
$ws.Range("B23").Value = "XPath"
$ws.Range("C23").Value = "//div[@id='ResetData']//a[text()='Close']"
